$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos Lab4")

# The student re-measured the "Insertion Sort" timings for n=1000 (row 2) and
# typed the new figures straight into B2:D2 (milliseconds, not the old toy
# Fibonacci-ish numbers). Typing plain numbers over cells that previously
# carried the "0.00" number style resets them back to the workbook's default
# (General) style, so clear the style before writing the new values.
$ws.Range("B2:D2").Style = "Normal"

$ws.Range("B2").Value = 609.38
$ws.Range("C2").Value = 640.63
$ws.Range("D2").Value = 31.25

# The dependent formulas in column C (C4:C11) and D (D3:D11) recalc
# automatically from the new seeds.

# Leave the selection on D2, matching where the student ended up after
# editing the row.
$ws.Range("D2").Select() | Out-Null
